# Insert a new weekly price record for "Cebollín" at row 188.
# All existing rows from 188 downward shift down by one (to 189..217).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 188, pushing the rest down.
$ws.Rows.Item(188).Insert()

# Populate the new row 188 with the new record's data.
$ws.Cells.Item(188, 1).Value  = 7
$ws.Cells.Item(188, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(188, 3).Value  = "Ñuble"
$ws.Cells.Item(188, 4).Value  = 45209
$ws.Cells.Item(188, 5).Value  = 16
$ws.Cells.Item(188, 6).Value  = 100112037
$ws.Cells.Item(188, 7).Value  = "Cebollín"
$ws.Cells.Item(188, 8).Value  = "Sin especificar"
$ws.Cells.Item(188, 9).Value  = "Primera"
$ws.Cells.Item(188, 10).Value = 150
$ws.Cells.Item(188, 11).Value = 6000
$ws.Cells.Item(188, 12).Value = 6000
$ws.Cells.Item(188, 13).Value = 6000
$ws.Cells.Item(188, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(188, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(188, 16).Value = 167
$ws.Cells.Item(188, 17).Value = 36
$ws.Cells.Item(188, 18).Value = "Hortaliza"
